$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.814.10'
$ws.Range('E2').Value = '  -1.79%  '
$ws.Range('D3').Value = '1.890.40'
$ws.Range('E3').Value = '  -1.53%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7682'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3121'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.65%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.31'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07213'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08074'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7650'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.498'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('D14').Value = '1.912.94'
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.25'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.146'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').Value = '29.839.31'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.92'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.27'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007756'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.03%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.158.18'
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.125'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1552'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.390'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.46'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.06%  '
$ws.Range('E28').Value = '  -2.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.039'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.437'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.550'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.465'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.100'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05499'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.96%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.259'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7459'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.624'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.48%  '
$ws.Range('E39').Value = '  -1.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.778'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('D41').Value = '1.137.62'
$ws.Range('E41').Value = '  +10.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '73.40'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4409'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.889'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8506'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.001'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.83'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.882'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.898'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.025'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +10.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.447'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.80%  '
